$d = $word.ActiveDocument

# Merge the three split runs ("Event organi" + "s" + "ers are prohibited...")
# into a single run by doing a Find/Replace over the full sentence.
$d.Content.Find.Execute(
    "Event organisers are prohibited from deleting events that already have registered participants.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Event organisers are prohibited from deleting events that already have registered participants.",
    2) | Out-Null

# Remove the "Users cannot unregister..." bullet paragraph and the trailing
# empty ListParagraph paragraph that followed it.
$target = "Users cannot unregister from an event once the registration deadline has passed."
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $text = $p.Range.Text.Trim()
    if ($text -eq $target) {
        # Delete this paragraph plus the very next paragraph (the now-empty
        # trailing ListParagraph) in one shot, including both paragraph marks.
        $start = $p.Range.Start
        $next = $p.Next()
        $end = $next.Range.End
        $d.Range($start, $end).Delete()
    }
}
